$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)

# --- 1. Give the "Content Placeholder 2" shape an explicit position/size ---
$sh.Left = 67.12496062992126
$sh.Top = 129.28111236220474
$sh.Width = 828.0
$sh.Height = 342.6250493700787

$tr = $sh.TextFrame.TextRange

# --- 2. Split " --user [username] --password [password] -" so that the two
#        credential placeholders are prefixed with an "Earthdata" run each ---
$full = $tr.Text
$runText = " --user [username] --password [password] -"
$runStart = $full.IndexOf($runText)
$runEnd = $runStart + $runText.Length

# 2a. " --user [" | "Earthdata username] --password [password] -"
$idx = $full.IndexOf("username")
$len1 = $runEnd - $idx
$sub = $tr.Characters($idx + 1, $len1)
$rest = $full.Substring($idx, $len1)
$sub.Text = "Earthdata " + $rest

# 2b. split off the new "Earthdata" run from " username] --password [password] -"
$full2 = $tr.Text
$runEnd2 = $runEnd + 10   # "Earthdata ".Length
$startAfterEarthdata = $idx + 9   # "Earthdata".Length
$len2 = $runEnd2 - $startAfterEarthdata
$sub2 = $tr.Characters($startAfterEarthdata + 1, $len2)
$rest2 = $full2.Substring($startAfterEarthdata, $len2)
$sub2.Text = $rest2

# 2c. locate the bracketed "[password]" (the second occurrence of the word)
$full3 = $tr.Text
$firstPass = $full3.IndexOf("password")
$secondPass = $full3.IndexOf("password", $firstPass + 1)
$len3 = $runEnd2 - $secondPass
$sub3 = $tr.Characters($secondPass + 1, $len3)
$rest3 = $full3.Substring($secondPass, $len3)
$sub3.Text = "Earthdata " + $rest3

# 2d. split off the second new "Earthdata" run from " password] -"
$full4 = $tr.Text
$runEnd3 = $runEnd2 + 10   # "Earthdata ".Length
$startAfterEarthdata2 = $secondPass + 9
$len4 = $runEnd3 - $startAfterEarthdata2
$sub4 = $tr.Characters($startAfterEarthdata2 + 1, $len4)
$rest4 = $full4.Substring($startAfterEarthdata2, $len4)
$sub4.Text = $rest4

# --- 3. Update the closing note to reference steps 5 and 6 instead of a script ---
$full5 = $tr.Text
$oldNote = "Note that this process will take days of uninterrupted internet connection. However, you can start the process where you left off using the script from step 5."
$newNote = "Note that this process will take days of uninterrupted internet connection. However, you can start the process where you left off using steps 5 and 6."
$noteStart = $full5.IndexOf($oldNote)
$subNote = $tr.Characters($noteStart + 1, $oldNote.Length)
$subNote.Text = $newNote
